# Facilitators guidelines - Moebius.docx : Swahili -> English text replacements
# plus default document language sw-KE -> sw-TZ.

$d = $word.ActiveDocument

# --- Table-cell / label text translations -------------------------------
# (old Swahili phrase -> new English phrase), each occurs exactly once
# except the last one, which repeats identically 6 times through the doc.
$replacements = @(
    @("Kichwa cha Video", "Video Title"),
    @("Mada", "Topic"),
    @("Jiometri", "Geometry"),
    @("Malengo", "Aim(s)"),
    @("Urefu", "Length"),
    @("Mahali pa Kambi", "Camp Location"),
    @("Wawezeshaji", "Facilitators"),
    @("N. ya wanafunzi", "N. of students"),
    @("Tarehe", "Date"),
    @("Rasilimali", "Resources"),
    @("inahitajika", "needed"),
    @("Maandalizi", "Preparations"),
    @("Muda wa video", "Video time"),
    @("Mwezeshaji anafanya nini", "What facilitator does"),
    @("Wanachofanya wanafunzi", "What learners do"),
    @("Utangulizi Mkuu wa Video ya VMC", "General VMC Video Introduction"),
    @("Utangulizi wa Video", "Video Introduction"),
    @("Utangulizi wa jaribio la kwanza", "Introduction of the first experiment"),
    @("Kusaidia mchakato, kuchochea mawazo", "Assist the process, provoke thoughts")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Default document language: sw-KE -> sw-TZ ---------------------------
$normal = $d.Styles.Item("Normal")
$normal.LanguageID = "sw-TZ"
